$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 4000
$ws.Range("J29").Value = 4000
$ws.Range("L29").Value = 12000
$ws.Range("N29").Value = -12562

# Row 39
$ws.Range("H38").Value = 88.666664
$ws.Range("J38").Value = 309
$ws.Range("L38").Value = 927
$ws.Range("N38").Value = -1671

# Row 43
$ws.Range("H43").Value = 499
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 499
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 499
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -637

# Row 58
$ws.Range("H58").Value = 600.8333
$ws.Range("I58").Value = 661
$ws.Range("J58").Value = 300
$ws.Range("K58").Value = 1983
$ws.Range("L58").Value = 900
$ws.Range("M58").Value = -1833
$ws.Range("N58").Value = -1200

# Row 113
$ws.Range("H113").Value = 38465620
$ws.Range("I113").Value = 52635068
$ws.Range("K113").Value = 52635068
$ws.Range("M113").Value = -52631814

# Row 129
$ws.Range("H129").Value = 213889.42
$ws.Range("J129").Value = 239312.92
$ws.Range("L129").Value = 717938.76
$ws.Range("N129").Value = -727938.76

# Row 132
$ws.Range("H132").Value = 3290.8928
$ws.Range("I132").Value = 3701.9565
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 11105.8695
$ws.Range("L132").Value = 4200
$ws.Range("M132").Value = -8575.869499999999
$ws.Range("N132").Value = -9260

# Row 137
$ws.Range("H137").Value = 2602.087
$ws.Range("I137").Value = 2235.8333
$ws.Range("J137").Value = 3920.6
$ws.Range("K137").Value = 6707.499899999999
$ws.Range("L137").Value = 11761.8
$ws.Range("M137").Value = -4157.499899999999
$ws.Range("N137").Value = -16861.8

# Row 138
$ws.Range("H138").Value = 14928375
$ws.Range("I138").Value = 50001708
$ws.Range("J138").Value = 3552
$ws.Range("K138").Value = 150005124
$ws.Range("L138").Value = 10656
$ws.Range("M138").Value = -149999984
$ws.Range("N138").Value = -20936

# Row 141
$ws.Range("H141").Value = 1338.3024
$ws.Range("I141").Value = 876.9729599999999
$ws.Range("K141").Value = 2630.91888
$ws.Range("M141").Value = 2549.08112

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1661.8572
$ws.Range("I2").Value = 1437.75
$ws.Range("K2").Value = 1437.75
$ws.Range("M2").Value = -1324.75

# Row 32
$ws.Range("H32").Value = 4407.829
$ws.Range("I32").Value = 2947.9866
$ws.Range("K32").Value = 2947.9866
$ws.Range("M32").Value = -2660.9866

# Row 45
$ws.Range("H45").Value = 2703.516
$ws.Range("I45").Value = 2689.842
$ws.Range("J45").Value = 2725.1667
$ws.Range("K45").Value = 2689.842
$ws.Range("L45").Value = 2725.1667
$ws.Range("M45").Value = -2312.842
$ws.Range("N45").Value = -3479.1667

# Row 61
$ws.Range("H61").Value = 347360.03
$ws.Range("I61").Value = 392265.4
$ws.Range("K61").Value = 392265.4
$ws.Range("M61").Value = -392053.4

# Row 74
$ws.Range("H74").Value = 28573352
$ws.Range("I74").Value = 30304944
$ws.Range("J74").Value = 2100
$ws.Range("K74").Value = 30304944
$ws.Range("L74").Value = 2100
$ws.Range("M74").Value = -30304070
$ws.Range("N74").Value = -3848

# Row 77
$ws.Range("H77").Value = 28573352
$ws.Range("I77").Value = 30304944
$ws.Range("J77").Value = 2100
$ws.Range("K77").Value = 151524720
$ws.Range("L77").Value = 10500
$ws.Range("M77").Value = -151520352
$ws.Range("N77").Value = -19236

# Row 116
$ws.Range("H116").Value = 1661.8572
$ws.Range("I116").Value = 1437.75
$ws.Range("K116").Value = 1437.75
$ws.Range("M116").Value = 856.25

# Row 122
$ws.Range("H122").Value = 2215.0908
$ws.Range("I122").Value = 1668
$ws.Range("K122").Value = 5004
$ws.Range("M122").Value = -2554

# Row 136
$ws.Range("H136").Value = 347360.03
$ws.Range("I136").Value = 392265.4
$ws.Range("K136").Value = 1176796.2
$ws.Range("M136").Value = -1174246.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1661.8572
$ws.Range("I3").Value = 1437.75
$ws.Range("K3").Value = 1437.75
$ws.Range("M3").Value = -1323.75

# Row 105
$ws.Range("H105").Value = 1821.3077
$ws.Range("I105").Value = 1793.7931
$ws.Range("J105").Value = 1901.1
$ws.Range("K105").Value = 1793.7931
$ws.Range("L105").Value = 1901.1
$ws.Range("M105").Value = -46.79310000000009
$ws.Range("N105").Value = -5395.1

# Row 107
$ws.Range("H107").Value = 630.1667
$ws.Range("I107").Value = 526.625
$ws.Range("J107").Value = 837.25
$ws.Range("K107").Value = 526.625
$ws.Range("L107").Value = 837.25
$ws.Range("M107").Value = 1393.375
$ws.Range("N107").Value = -4677.25

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 130
$ws.Range("H130").Value = 53555
$ws.Range("J130").Value = 53555
$ws.Range("L130").Value = 53555
$ws.Range("N130").Value = -63595

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4201.9
$ws.Range("I31").Value = 2777.24
$ws.Range("J31").Value = 6576.3335
$ws.Range("K31").Value = 2777.24
$ws.Range("L31").Value = 6576.3335
$ws.Range("M31").Value = -2482.24
$ws.Range("N31").Value = -7166.3335

# Row 34
$ws.Range("H34").Value = 4201.9
$ws.Range("I34").Value = 2777.24
$ws.Range("J34").Value = 6576.3335
$ws.Range("K34").Value = 2777.24
$ws.Range("L34").Value = 6576.3335
$ws.Range("M34").Value = -2575.24
$ws.Range("N34").Value = -6980.3335

# Row 58
$ws.Range("H58").Value = 10056.891
$ws.Range("I58").Value = 843.0789
$ws.Range("J58").Value = 30652.47
$ws.Range("K58").Value = 843.0789
$ws.Range("L58").Value = 30652.47
$ws.Range("M58").Value = -640.0789
$ws.Range("N58").Value = -31058.47

# Row 122
$ws.Range("H122").Value = 5500.3335
$ws.Range("I122").Value = 6500.5
$ws.Range("K122").Value = 19501.5
$ws.Range("M122").Value = -17051.5

# Row 132
$ws.Range("H132").Value = 1966.7646
$ws.Range("I132").Value = 1397.6957
$ws.Range("K132").Value = 4193.0871
$ws.Range("M132").Value = -1663.0871

# Row 134
$ws.Range("H134").Value = 1015.1539
$ws.Range("I134").Value = 938.26086
$ws.Range("K134").Value = 2814.78258
$ws.Range("M134").Value = -279.7825800000001

# Row 136
$ws.Range("H136").Value = 10056.891
$ws.Range("I136").Value = 843.0789
$ws.Range("J136").Value = 30652.47
$ws.Range("K136").Value = 2529.2367
$ws.Range("L136").Value = 91957.41
$ws.Range("M136").Value = 20.76330000000007
$ws.Range("N136").Value = -97057.41

$ws = $wb.Worksheets.Item("CUL")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# Row 126
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = 440

# Row 129
$ws.Range("H129").Value = 239271.53
$ws.Range("I129").Value = 988
$ws.Range("J129").Value = 455892.9
$ws.Range("K129").Value = 2964
$ws.Range("L129").Value = 1367678.7
$ws.Range("M129").Value = 2036
$ws.Range("N129").Value = -1377678.7

# Row 131
$ws.Range("H131").Value = 755.76
$ws.Range("J131").Value = 767.1158
$ws.Range("L131").Value = 2301.3474
$ws.Range("N131").Value = -12381.3474

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 19232666
$ws.Range("I102").Value = 20834908
$ws.Range("K102").Value = 20834908
$ws.Range("M102").Value = -20833286

# Row 122
$ws.Range("H122").Value = 111113496
$ws.Range("I122").Value = 37038884
$ws.Range("K122").Value = 111116652
$ws.Range("M122").Value = -111114202

# Row 126
$ws.Range("H126").Value = 4736.1333
$ws.Range("J126").Value = 5598.857
$ws.Range("L126").Value = 16796.571
$ws.Range("N126").Value = -21736.571

# Row 132
$ws.Range("H132").Value = 17949.818
$ws.Range("I132").Value = 3378.4092
$ws.Range("J132").Value = 47092.637
$ws.Range("K132").Value = 10135.2276
$ws.Range("L132").Value = 141277.911
$ws.Range("M132").Value = -7605.2276
$ws.Range("N132").Value = -146337.911

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2936.7812
$ws.Range("I40").Value = 2710.2222
$ws.Range("J40").Value = 4160.2
$ws.Range("K40").Value = 2710.2222
$ws.Range("L40").Value = 4160.2
$ws.Range("M40").Value = -2574.2222
$ws.Range("N40").Value = -4432.2

# Row 132
$ws.Range("H132").Value = 1186.9048
$ws.Range("I132").Value = 1079.54
$ws.Range("K132").Value = 3238.62
$ws.Range("M132").Value = -708.6199999999999

# Row 136
$ws.Range("H136").Value = 973.7308
$ws.Range("I136").Value = 973.7308
$ws.Range("K136").Value = 2921.1924
$ws.Range("M136").Value = -371.1923999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 719.31885
$ws.Range("I132").Value = 485
$ws.Range("J132").Value = 1954.8182
$ws.Range("K132").Value = 1455
$ws.Range("L132").Value = 5864.4546
$ws.Range("M132").Value = 1075
$ws.Range("N132").Value = -10924.4546

# Row 136
$ws.Range("H136").Value = 14494540
$ws.Range("I136").Value = 22728166
$ws.Range("J136").Value = 3356.6
$ws.Range("K136").Value = 68184498
$ws.Range("L136").Value = 10069.8
$ws.Range("M136").Value = -68181948
$ws.Range("N136").Value = -15169.8
